# aggiornamento a 9/09 compreso
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(367, 44441, 0, 2, 61.06870229007634),
    @(368, 44442, 2, 4, 122.1374045801527),
    @(369, 44443, 2, 6, 183.206106870229),
    @(370, 44444, 0, 6, 183.206106870229),
    @(371, 44445, 1, 5, 152.6717557251908),
    @(372, 44446, 0, 5, 152.6717557251908),
    @(373, 44447, 1, 6, 183.206106870229),
    @(374, 44448, 0, 6, 183.206106870229)
)

foreach ($entry in $data) {
    $r = $entry[0]

    # Carry the date column's formatting (style index 2) down from the
    # previous row, same as dragging the fill handle in the UI.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

$excel.CutCopyMode = 0
